# Adds two new columns, I ("I0") and J ("IF"), to the sheet, with header
# cells styled like the rest of the header row, and fills in the per-row
# values for rows 2 through 62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers "I0" and "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font + thin border + centered alignment) from
# the existing header style (column H) onto the two new header cells so
# they share the same cell style as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows 2-62: fill in the I and J values ---
$data = @(
    @(2,3,4),
    @(3,11,11),
    @(4,5,5),
    @(5,6,7),
    @(6,8,8),
    @(7,6,6),
    @(8,6,7),
    @(9,8,8),
    @(10,6,7),
    @(11,6,7),
    @(12,7,8),
    @(13,9,9),
    @(14,8,8),
    @(15,6,7),
    @(16,6,8),
    @(17,7,7),
    @(18,5,6),
    @(19,7,7),
    @(20,8,8),
    @(21,7,8),
    @(22,6,7),
    @(23,6,6),
    @(24,6,8),
    @(25,6,7),
    @(26,8,8),
    @(27,11,11),
    @(28,8,8),
    @(29,7,7),
    @(30,7,7),
    @(31,10,10),
    @(32,7,8),
    @(33,6,6),
    @(34,4,6),
    @(35,7,7),
    @(36,6,6),
    @(37,8,8),
    @(38,4,5),
    @(39,6,6),
    @(40,1,3),
    @(41,8,8),
    @(42,5,5),
    @(43,5,6),
    @(44,8,8),
    @(45,8,8),
    @(46,8,8),
    @(47,6,6),
    @(48,9,9),
    @(49,7,7),
    @(50,6,6),
    @(51,8,8),
    @(52,3,4),
    @(53,6,6),
    @(54,11,11),
    @(55,5,5),
    @(56,9,9),
    @(57,4,5),
    @(58,4,4),
    @(59,8,8),
    @(60,9,9),
    @(61,5,6),
    @(62,8,8)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
